$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume cells are treated as text so values like "147.00" or
# "0.650" keep their exact formatting instead of being parsed as numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.671.84'
$ws.Range('E2').Value = '  -2.28%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.591.70'
$ws.Range('E3').Value = '  -2.85%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.01'
$ws.Range('E5').Value = '  -2.69%  '
$ws.Range('E6').Value = '  -2.40%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('E8').Value = '  -3.09%  '
$ws.Range('E9').Value = '  -1.81%  '
$ws.Range('E10').Value = '  -3.90%  '
$ws.Range('E11').Value = '  -1.77%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.813.74'
$ws.Range('E12').Value = '  -2.79%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.591.45'
$ws.Range('E13').Value = '  -2.88%  '
$ws.Range('E14').Value = '  -2.94%  '
$ws.Range('E15').Value = '  -4.05%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.67'
$ws.Range('E16').Value = '  -1.25%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.647.56'
$ws.Range('E17').Value = '  -2.15%  '
$ws.Range('E18').Value = '  -1.99%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '208.01'
$ws.Range('E19').Value = '  -5.10%  '
$ws.Range('E20').Value = '  +0.01%  '
$ws.Range('E21').Value = '  -3.89%  '
$ws.Range('E22').Value = '  -3.59%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.39'
$ws.Range('E23').Value = '  -1.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.92'
$ws.Range('E24').Value = '  -2.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.00'
$ws.Range('E25').Value = '  -0.57%  '
$ws.Range('E26').Value = '  +0.27%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.29'
$ws.Range('E27').Value = '  -0.61%  '
$ws.Range('E28').Value = '  -3.83%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.30'
$ws.Range('E29').Value = '  -2.52%  '
$ws.Range('E30').Value = '  -0.92%  '
$ws.Range('E31').Value = '  -2.39%  '
$ws.Range('E32').Value = '  -4.66%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.650'
$ws.Range('E33').Value = '  +18.48%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.90'
$ws.Range('E34').Value = '  -3.61%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.314.88'
$ws.Range('E35').Value = '  -1.74%  '
$ws.Range('E36').Value = '  -4.98%  '
$ws.Range('E37').Value = '  -1.83%  '
$ws.Range('E38').Value = '  -2.59%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.830'
$ws.Range('E39').Value = '  -2.86%  '
$ws.Range('E40').Value = '  +0.08%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.39'
$ws.Range('E41').Value = '  +2.39%  '
$ws.Range('E42').Value = '  -2.08%  '
$ws.Range('E43').Value = '  -3.96%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.48'
$ws.Range('E44').Value = '  -1.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.726.96'
$ws.Range('E45').Value = '  -2.57%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '89.98'
$ws.Range('E46').Value = '  -1.19%  '
$ws.Range('E47').Value = '  -1.53%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.840'
$ws.Range('E48').Value = '  +3.69%  '
$ws.Range('E49').Value = '  -0.41%  '
$ws.Range('E50').Value = '  -1.90%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.53'
$ws.Range('E51').Value = '  -1.30%  '
